$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the paragraph that holds the "energetic/dynamic ... (full of
# beans) ... slow/indolent" line (the one the _GoBack bookmark sits in).
# ------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "full of beans") {
        $target = $p
        break
    }
}

# ------------------------------------------------------------------
# Step 1: move the "_GoBack" bookmark so it sits at the very end of the
# paragraph's text (right after "indolent") instead of right after
# "(full of beans)   ".
#
# A collapsed Range placed exactly at "end of paragraph - 1" (i.e. right
# before the paragraph mark) trips a bug in Bookmarks.Add on this host,
# so we can't just do $d.Bookmarks.Add("_GoBack", $collapsedRange) there
# directly. Work around it by inserting a one-character placeholder at
# that spot (which safely shifts the paragraph end outward), wrapping
# the bookmark around the placeholder, and then clearing the
# placeholder's text from inside the bookmark's own Range so it
# collapses back down to a zero-length bookmark in the right place -
# exactly like Word does when bookmarked text is deleted from within.
# ------------------------------------------------------------------
$find = $target.Range.Duplicate
$find.Find.Execute("indolent", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$find.Collapse(0)
$find.InsertAfter("X")

$d.Bookmarks.Add("_GoBack", $find)
$bm = $d.Bookmarks.Item("_GoBack")
$bmRange = $bm.Range
$bmRange.Text = ""

# ------------------------------------------------------------------
# Step 2: delete the trailing empty paragraph that followed this one by
# removing this paragraph's own end-of-paragraph mark, which merges it
# with the (empty) paragraph after it.
# ------------------------------------------------------------------
$mark = $d.Range($target.Range.End - 1, $target.Range.End)
$mark.Delete()
